$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

# Row 3: primary phone number placeholder -> realistic UK number
$ws.Range("H3").Value = "+44 1632 960852"

# Row 4: primary + alternative phone numbers -> realistic numbers
$ws.Range("H4").Value = "+1-613-555-0182"
$ws.Range("I4").Value = "+36 55 979 922"

# Rows 5-29: the per-row placeholder phone numbers (888-888-88XX) are
# replaced by re-using the two realistic numbers set on rows 3/4,
# alternating odd/even the same way the original placeholders did.
for ($r = 5; $r -le 29; $r++) {
    if ($r % 2 -eq 1) {
        $ws.Range("H$r").Value = "+44 1632 960852"
    } else {
        $ws.Range("H$r").Value = "+1-613-555-0182"
    }
}
